$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1780
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 1733.3334
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 1733.3334
$ws.Range("M40").Value = -1625
$ws.Range("N40").Value = -2083.3334
$ws.Range("H64").Value = 3386.0833
$ws.Range("I64").Value = 2932
$ws.Range("J64").Value = 3710.4285
$ws.Range("K64").Value = 2932
$ws.Range("L64").Value = 3710.4285
$ws.Range("M64").Value = -2684
$ws.Range("N64").Value = -4206.4285
$ws.Range("H67").Value = 3386.0833
$ws.Range("I67").Value = 2932
$ws.Range("J67").Value = 3710.4285
$ws.Range("K67").Value = 2932
$ws.Range("L67").Value = 3710.4285
$ws.Range("M67").Value = -2074
$ws.Range("N67").Value = -5426.4285
$ws.Range("H116").Value = 16787.143
$ws.Range("I116").Value = 26127.5
$ws.Range("J116").Value = 4333.3335
$ws.Range("K116").Value = 26127.5
$ws.Range("L116").Value = 4333.3335
$ws.Range("M116").Value = -22685.5
$ws.Range("N116").Value = -11217.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 3000
$ws.Range("J11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3288
$ws.Range("H34").Value = 9000
$ws.Range("J34").Value = 9000
$ws.Range("L34").Value = 9000
$ws.Range("N34").Value = -9542
$ws.Range("H63").Value = 5790.875
$ws.Range("I63").Value = 4636.2856
$ws.Range("J63").Value = 6688.8887
$ws.Range("K63").Value = 4636.2856
$ws.Range("L63").Value = 6688.8887
$ws.Range("M63").Value = -3950.2856
$ws.Range("N63").Value = -8060.8887
$ws.Range("H66").Value = 5790.875
$ws.Range("I66").Value = 4636.2856
$ws.Range("J66").Value = 6688.8887
$ws.Range("K66").Value = 23181.428
$ws.Range("L66").Value = 33444.4435
$ws.Range("M66").Value = -19749.428
$ws.Range("N66").Value = -40308.4435
$ws.Range("H74").Value = 1251.2142
$ws.Range("I74").Value = 1178.909
$ws.Range("K74").Value = 1178.909
$ws.Range("M74").Value = -304.9090000000001
$ws.Range("H77").Value = 1251.2142
$ws.Range("I77").Value = 1178.909
$ws.Range("K77").Value = 5894.545
$ws.Range("M77").Value = -1526.545
$ws.Range("H82").Value = 14658
$ws.Range("I82").Value = 2164
$ws.Range("J82").Value = 20905
$ws.Range("K82").Value = 2164
$ws.Range("L82").Value = 20905
$ws.Range("M82").Value = -1803
$ws.Range("N82").Value = -21627
$ws.Range("H85").Value = 14658
$ws.Range("I85").Value = 2164
$ws.Range("J85").Value = 20905
$ws.Range("K85").Value = 2164
$ws.Range("L85").Value = 20905
$ws.Range("M85").Value = -916
$ws.Range("N85").Value = -23401
$ws.Range("H132").Value = 16411309
$ws.Range("I132").Value = 20001284
$ws.Range("J132").Value = 93239.82000000001
$ws.Range("K132").Value = 60003852
$ws.Range("L132").Value = 279719.46
$ws.Range("M132").Value = -60001322
$ws.Range("N132").Value = -284779.46

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9282.909
$ws.Range("I134").Value = 3160.8823
$ws.Range("J134").Value = 30097.8
$ws.Range("K134").Value = 9482.6469
$ws.Range("L134").Value = 90293.39999999999
$ws.Range("M134").Value = -6947.6469
$ws.Range("N134").Value = -95363.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2225
$ws.Range("I62").Value = 1300
$ws.Range("J62").Value = 2533.3333
$ws.Range("K62").Value = 1300
$ws.Range("L62").Value = 2533.3333
$ws.Range("M62").Value = -676
$ws.Range("N62").Value = -3781.3333
$ws.Range("H65").Value = 2225
$ws.Range("I65").Value = 1300
$ws.Range("J65").Value = 2533.3333
$ws.Range("K65").Value = 6500
$ws.Range("L65").Value = 12666.6665
$ws.Range("M65").Value = -3380
$ws.Range("N65").Value = -18906.6665
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -22246
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -71232
$ws.Range("H122").Value = 917.04346
$ws.Range("I122").Value = 556.5714
$ws.Range("J122").Value = 1477.7778
$ws.Range("K122").Value = 1669.7142
$ws.Range("L122").Value = 4433.3334
$ws.Range("M122").Value = 780.2857999999999
$ws.Range("N122").Value = -9333.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 5968.4614
$ws.Range("J102").Value = 5968.4614
$ws.Range("L102").Value = 17905.3842
$ws.Range("N102").Value = -22773.3842
$ws.Range("H117").Value = 1192
$ws.Range("I117").Value = 976.3333
$ws.Range("J117").Value = 1250.8182
$ws.Range("K117").Value = 2928.9999
$ws.Range("L117").Value = 3752.4546
$ws.Range("M117").Value = 513.0001000000002
$ws.Range("N117").Value = -10636.4546
$ws.Range("H129").Value = 23811492
$ws.Range("I129").Value = 2415
$ws.Range("J129").Value = 27779670
$ws.Range("K129").Value = 7245
$ws.Range("L129").Value = 83339010
$ws.Range("M129").Value = -2245
$ws.Range("N129").Value = -83349010

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 5000627
$ws.Range("J24").Value = 1253.5
$ws.Range("L24").Value = 1253.5
$ws.Range("N24").Value = -1599.5
$ws.Range("H45").Value = 19999.8
$ws.Range("J45").Value = 19999.8
$ws.Range("L45").Value = 19999.8
$ws.Range("N45").Value = -21117.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2872.0527
$ws.Range("I7").Value = 1371.2858
$ws.Range("J7").Value = 3747.5
$ws.Range("K7").Value = 1371.2858
$ws.Range("L7").Value = 3747.5
$ws.Range("M7").Value = -1259.2858
$ws.Range("N7").Value = -3971.5
$ws.Range("H20").Value = 2999.6667
$ws.Range("J20").Value = 2999.6667
$ws.Range("L20").Value = 2999.6667
$ws.Range("N20").Value = -3451.6667
$ws.Range("H100").Value = 3551
$ws.Range("I100").Value = 1400
$ws.Range("J100").Value = 3858.2856
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 3858.2856
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -4940.2856
$ws.Range("H126").Value = 2872.0527
$ws.Range("I126").Value = 1371.2858
$ws.Range("J126").Value = 3747.5
$ws.Range("K126").Value = 4113.857400000001
$ws.Range("L126").Value = 11242.5
$ws.Range("M126").Value = -1643.857400000001
$ws.Range("N126").Value = -16182.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 11468
$ws.Range("J13").Value = 7202
$ws.Range("L13").Value = 7202
$ws.Range("N13").Value = -7482
$ws.Range("H64").Value = 9000
$ws.Range("J64").Value = 9000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9496
$ws.Range("H67").Value = 9000
$ws.Range("J67").Value = 9000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10716
$ws.Range("H130").Value = 29933
$ws.Range("J130").Value = 29933
$ws.Range("L130").Value = 29933
$ws.Range("N130").Value = -39973
$ws.Range("H132").Value = 52098216
$ws.Range("I132").Value = 64572616
$ws.Range("J132").Value = 3586663
$ws.Range("K132").Value = 193717848
$ws.Range("L132").Value = 10759989
$ws.Range("M132").Value = -193715318
$ws.Range("N132").Value = -10765049